$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tracking")

$rng = $ws.Range("A2:F27")
$key = $ws.Range("E2:E27")

$rng.Sort($key, 1)

$ws.Range("G21:G27").ClearContents()
